$wb = $excel.ActiveWorkbook

$msg = "Validado com sucesso! Nenhuma divergência entre o SPED e o relatório foi encontrada!"

# ----- Sheet "Bico" (first sheet) -----
$wsBico = $wb.Worksheets.Item(1)

$wsBico.Cells.Item(1, 8).Value = "Obs_relatorio"
$wsBico.Cells.Item(1, 9).Value = "Obs_sped"

for ($r = 2; $r -le 12; $r++) {
    $wsBico.Cells.Item($r, 8).Value = $msg
    $wsBico.Cells.Item($r, 9).Font.Bold = $false
}

# ----- Sheet "Tanque" (second sheet) -----
$wsTanque = $wb.Worksheets.Item(2)

$wsTanque.Cells.Item(1, 6).Value = "Obs_relatorio"
$wsTanque.Cells.Item(1, 7).Value = "Obs_sped"

for ($r = 2; $r -le 7; $r++) {
    $wsTanque.Cells.Item($r, 6).Value = $msg
    $wsTanque.Cells.Item($r, 7).Font.Bold = $false
}
